$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value of F3 to append ",a,b" info to the path (new text)
$ws.Range("F3").Value = "C:\temp\a,b\"

# Select F3 to match the new active cell / selection in the sheet view
$ws.Range("F3").Select()
